$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily log. It belongs right
# before the existing row that currently sits at row 166 (chronologically /
# as captured by the source system), so insert a fresh row there and push
# every following record down by one (166->167, ..., 200->201).
$ws.Rows.Item(166).Insert()

# Fill in the newly inserted row with the new observation.
$ws.Range("A166").Value = 9
$ws.Range("B166").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C166").Value = "Metropolitana"
$ws.Range("D166").Value = 44663
$ws.Range("E166").Value = 13
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100101
$ws.Range("H166").Value = "Berries"
$ws.Range("I166").Value = 100101001
$ws.Range("J166").Value = "Arándano (blue)"
$ws.Range("K166").Value = "Sin especificar"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 580
$ws.Range("N166").Value = 4500
$ws.Range("O166").Value = 5000
$ws.Range("P166").Value = 4741
$ws.Range("Q166").Value = "$/bandeja 2 kilos"
$ws.Range("R166").Value = "Provincia de Linares"
$ws.Range("S166").Value = 2370
$ws.Range("T166").Value = 2
